{"js": "// Cleanup pass: remove decorative \"picture\" paragraphs, the faint dashed\n// separator rules, and the stray empty spacer paragraphs (spacing-before=40\n// twips, typically left sitting right after a table) that accumulated in\n// this document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].inlinePictures.load(\"items\");\n}\nawait context.sync();\n\nconst SEP_CHAR = \"\\u2500\"; // \u2500 box-drawing light horizontal, used as a visual separator rule\nconst toDelete = [];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text;\n  const hasPicture = p.inlinePictures.items.length > 0;\n\n  const isImageParagraph = hasPicture;\n  const isSeparatorParagraph = text.length > 0 && text.split(SEP_CHAR).join(\"\") === \"\";\n  const isEmptySpacerParagraph = text === \"\" && !hasPicture;\n\n  if (isImageParagraph || isSeparatorParagraph || isEmptySpacerParagraph) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Cleanup pass: remove decorative \"picture\" paragraphs, the faint dashed\n# separator rules, and the stray empty spacer paragraphs (spacing-before=40\n# twips, typically left sitting right after a table) that accumulated in\n# this document.\n\n$d = $word.ActiveDocument\n$sepChar = [char]0x2500\n$cr = [char]13\n\n$count = $d.Paragraphs.Count\n$toDelete = New-Object System.Collections.ArrayList\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $range = $p.Range\n    $text = $range.Text\n    $shapeCount = $range.InlineShapes.Count\n\n    $isImageParagraph = $shapeCount -gt 0\n    $isSeparatorParagraph = ($text.Replace($sepChar, \"\")) -eq $cr\n    $isEmptySpacerParagraph = ($text -eq $cr) -and ($shapeCount -eq 0)\n\n    if ($isImageParagraph -or $isSeparatorParagraph -or $isEmptySpacerParagraph) {\n        [void]$toDelete.Add($i)\n    }\n}\n\n# Delete from the last index to the first so earlier indices stay valid.\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $idx = $toDelete[$j]\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
